$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.58"
$ws.Range("E2").Value = "'0.94%"
$ws.Range("D3").Value = "'35.84"
$ws.Range("E3").Value = "'1.45%"
$ws.Range("D4").Value = "'5.066"
$ws.Range("E4").Value = "'0.09%"
$ws.Range("D5").Value = "'0.08059"
$ws.Range("E5").Value = "'1.62%"
$ws.Range("D6").Value = "'1.943"
$ws.Range("E6").Value = "'2.55%"
$ws.Range("D7").Value = "'4.156"
$ws.Range("E7").Value = "'2.74%"
$ws.Range("D8").Value = "'7.834"
$ws.Range("E8").Value = "'0.78%"
$ws.Range("D9").Value = "'0.9291"
$ws.Range("E9").Value = "'0.17%"
$ws.Range("D10").Value = "'0.1310"
$ws.Range("E10").Value = "'-5.39%"
$ws.Range("D11").Value = "'0.1905"
$ws.Range("E11").Value = "'0.02%"
$ws.Range("D12").Value = "'0.09197"
$ws.Range("E12").Value = "'0.44%"
$ws.Range("D13").Value = "'0.03477"
$ws.Range("E13").Value = "'1.32%"
$ws.Range("D14").Value = "'0.09907"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("D15").Value = "'0.001413"
$ws.Range("E15").Value = "'1.76%"
$ws.Range("D16").Value = "'0.006649"
$ws.Range("E16").Value = "'14.03%"
$ws.Range("D17").Value = "'3.611"
$ws.Range("E18").Value = "'3.02%"
$ws.Range("E19").Value = "'0.48%"
$ws.Range("E20").Value = "'3.27%"
$ws.Range("D21").Value = "'5.163"
$ws.Range("E21").Value = "'2.38%"
$ws.Range("D22").Value = "'0.2532"
$ws.Range("E22").Value = "'5.63%"
$ws.Range("D23").Value = "'0.04413"
$ws.Range("E23").Value = "'-2.20%"
$ws.Range("D24").Value = "'0.001237"
$ws.Range("E24").Value = "'2.02%"
$ws.Range("D25").Value = "'0.004699"
$ws.Range("E25").Value = "'-1.51%"
$ws.Range("E26").Value = "'5.76%"
$ws.Range("E27").Value = "'4.35%"
$ws.Range("D39").Value = "'0.01995"
$ws.Range("E39").Value = "'7.06%"
$ws.Range("D40").Value = "'0.05209"
$ws.Range("E40").Value = "'9.11%"
$ws.Range("D41").Value = "'0.007615"
$ws.Range("E41").Value = "'3.86%"
$ws.Range("D42").Value = "'0.01013"
$ws.Range("E42").Value = "'5.16%"
$ws.Range("D43").Value = "'0.1364"
$ws.Range("E43").Value = "'3.02%"
$ws.Range("D44").Value = "'0.002102"
$ws.Range("D45").Value = "'0.01074"
$ws.Range("E45").Value = "'-2.42%"
$ws.Range("D46").Value = "'0.00006294"
$ws.Range("E46").Value = "'0.77%"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("E48").Value = "'0.46%"
$ws.Range("D49").Value = "'0.001602"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E51").Value = "'0.06%"
